$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update version, date, publisher/jurisdiction info ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date updated
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now filled in
$ws1.Range("B9").Value = "Alvearie Team"

# Row that used to be "Contact" / "No display for ContactDetail" becomes
# "Jurisdiction" / "United States of America"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# The old duplicate "Contact" row (row 11) is removed entirely, shifting
# everything below it up by one row.
$ws1.Rows("11:11").Delete()

# --- Sheet "Elements": update Short/Definition text for the root Extension row ---
$ws2 = $wb.Worksheets.Item("Elements")

$ws2.Range("K2").Value = "Days To Process"
$ws2.Range("L2").Value = "Number of days between the date the claim was received and the date the claim was processed"
